$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALERTS")
$ws.Range("A6").Value = "'2026-01-28"
$ws.Range("A6").Style = "Normal"
$ws.Range("B6").Value = "17:49:35"
$ws.Range("C6").Value = "17:00"
$ws.Range("D6").Value = "Bathroom"
$ws.Range("E6").Value = "MINIMAL"
$ws.Range("F6").Value = "MINIMAL ALERT: Bathroom occupied, no motion > 20s."

$ws.Range("A7").Value = "'2026-01-28"
$ws.Range("A7").Style = "Normal"
$ws.Range("B7").Value = "17:49:51"
$ws.Range("C7").Value = "17:00"
$ws.Range("D7").Value = "Bathroom"
$ws.Range("E7").Value = "MODERATE"
$ws.Range("F7").Value = "MODERATE ALERT: Bathroom occupied, no motion > 40s."

$ws.Range("A8").Value = "'2026-01-28"
$ws.Range("A8").Style = "Normal"
$ws.Range("B8").Value = "17:50:10"
$ws.Range("C8").Value = "17:00"
$ws.Range("D8").Value = "Bathroom"
$ws.Range("E8").Value = "CRITICAL"
$ws.Range("F8").Value = "CRITICAL ALERT: Bathroom occupied, no motion > 60s."

$ws = $wb.Worksheets.Item("Proximity")
$ws.Range("A14").Value = "'2026-01-28"
$ws.Range("A14").Style = "Normal"
$ws.Range("B14").Value = "17:49:48"
$ws.Range("C14").Value = "17:00"
$ws.Range("D14").Value = "Living Room Main Door"
$ws.Range("E14").Value = "ENTER"
$ws.Range("F14").Value = "User ENTERED Living Room Main Door"

$ws.Range("A15").Value = "'2026-01-28"
$ws.Range("A15").Style = "Normal"
$ws.Range("B15").Value = "17:49:50"
$ws.Range("C15").Value = "17:00"
$ws.Range("D15").Value = "Living Room Main Door"
$ws.Range("E15").Value = "EXIT"
$ws.Range("F15").Value = "User EXITED Living Room Main Door"

$ws.Range("A16").Value = "'2026-01-28"
$ws.Range("A16").Style = "Normal"
$ws.Range("B16").Value = "17:50:00"
$ws.Range("C16").Value = "17:00"
$ws.Range("D16").Value = "Living Room Main Door"
$ws.Range("E16").Value = "ENTER"
$ws.Range("F16").Value = "User ENTERED Living Room Main Door"

$ws = $wb.Worksheets.Item("Camera")
$ws.Range("A9").Value = "'2026-01-28"
$ws.Range("A9").Style = "Normal"
$ws.Range("B9").Value = "17:49:50"
$ws.Range("C9").Value = "17:00"
$ws.Range("D9").Value = "Living Room Main Door"
$ws.Range("E9").Value = "Image Captured"
$ws.Range("F9").Value = "Active"

$ws.Range("A10").Value = "'2026-01-28"
$ws.Range("A10").Style = "Normal"
$ws.Range("B10").Value = "17:50:02"
$ws.Range("C10").Value = "17:00"
$ws.Range("D10").Value = "Living Room Main Door"
$ws.Range("E10").Value = "Image Captured"
$ws.Range("F10").Value = "Active"

$ws = $wb.Worksheets.Item("Sleep")
$ws.Range("A10").Value = "'2026-01-28"
$ws.Range("A10").Style = "Normal"
$ws.Range("B10").Value = "17:49:35"
$ws.Range("C10").Value = "17:00"
$ws.Range("D10").Value = "Bedroom"
$ws.Range("E10").Value = "In Bed"
$ws.Range("F10").Value = 71
$ws.Range("G10").Value = 23
$ws.Range("H10").Value = "Occupied"

$ws.Range("A11").Value = "'2026-01-28"
$ws.Range("A11").Style = "Normal"
$ws.Range("B11").Value = "17:49:36"
$ws.Range("C11").Value = "17:00"
$ws.Range("D11").Value = "Bedroom"
$ws.Range("E11").Value = "In Bed"
$ws.Range("F11").Value = 50
$ws.Range("G11").Value = 2
$ws.Range("H11").Value = "Occupied"

$ws.Range("A12").Value = "'2026-01-28"
$ws.Range("A12").Style = "Normal"
$ws.Range("B12").Value = "17:49:36"
$ws.Range("C12").Value = "17:00"
$ws.Range("D12").Value = "Bedroom"
$ws.Range("E12").Value = "In Bed"
$ws.Range("F12").Value = 55
$ws.Range("G12").Value = 7
$ws.Range("H12").Value = "Occupied"

$ws.Range("A13").Value = "'2026-01-28"
$ws.Range("A13").Style = "Normal"
$ws.Range("B13").Value = "17:49:36"
$ws.Range("C13").Value = "17:00"
$ws.Range("D13").Value = "Bedroom"
$ws.Range("E13").Value = "In Bed"
$ws.Range("F13").Value = 50
$ws.Range("G13").Value = 2
$ws.Range("H13").Value = "Occupied"

$ws.Range("A14").Value = "'2026-01-28"
$ws.Range("A14").Style = "Normal"
$ws.Range("B14").Value = "17:49:37"
$ws.Range("C14").Value = "17:00"
$ws.Range("D14").Value = "Bedroom"
$ws.Range("E14").Value = "In Bed"
$ws.Range("F14").Value = 55
$ws.Range("G14").Value = 7
$ws.Range("H14").Value = "Occupied"

$ws.Range("A15").Value = "'2026-01-28"
$ws.Range("A15").Style = "Normal"
$ws.Range("B15").Value = "17:49:38"
$ws.Range("C15").Value = "17:00"
$ws.Range("D15").Value = "Bedroom"
$ws.Range("E15").Value = "In Bed"
$ws.Range("F15").Value = 50
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = "Occupied"

$ws.Range("A16").Value = "'2026-01-28"
$ws.Range("A16").Style = "Normal"
$ws.Range("B16").Value = "17:49:41"
$ws.Range("C16").Value = "17:00"
$ws.Range("D16").Value = "Bedroom"
$ws.Range("E16").Value = "In Bed"
$ws.Range("F16").Value = 49
$ws.Range("G16").Value = 1
$ws.Range("H16").Value = "Occupied"

$ws.Range("A17").Value = "'2026-01-28"
$ws.Range("A17").Style = "Normal"
$ws.Range("B17").Value = "17:49:43"
$ws.Range("C17").Value = "17:00"
$ws.Range("D17").Value = "Bedroom"
$ws.Range("E17").Value = "In Bed"
$ws.Range("F17").Value = 50
$ws.Range("G17").Value = 2
$ws.Range("H17").Value = "Occupied"

$ws.Range("A18").Value = "'2026-01-28"
$ws.Range("A18").Style = "Normal"
$ws.Range("B18").Value = "17:49:52"
$ws.Range("C18").Value = "17:00"
$ws.Range("D18").Value = "Bedroom"
$ws.Range("E18").Value = "In Bed"
$ws.Range("F18").Value = 56
$ws.Range("G18").Value = 8
$ws.Range("H18").Value = "Occupied"

$ws.Range("A19").Value = "'2026-01-28"
$ws.Range("A19").Style = "Normal"
$ws.Range("B19").Value = "17:49:54"
$ws.Range("C19").Value = "17:00"
$ws.Range("D19").Value = "Bedroom"
$ws.Range("E19").Value = "In Bed"
$ws.Range("F19").Value = 50
$ws.Range("G19").Value = 2
$ws.Range("H19").Value = "Occupied"

$ws.Range("A20").Value = "'2026-01-28"
$ws.Range("A20").Style = "Normal"
$ws.Range("B20").Value = "17:50:07"
$ws.Range("C20").Value = "17:00"
$ws.Range("D20").Value = "Bedroom"
$ws.Range("E20").Value = "In Bed"
$ws.Range("F20").Value = 49
$ws.Range("G20").Value = 1
$ws.Range("H20").Value = "Occupied"

$ws.Range("A21").Value = "'2026-01-28"
$ws.Range("A21").Style = "Normal"
$ws.Range("B21").Value = "17:50:24"
$ws.Range("C21").Value = "17:00"
$ws.Range("D21").Value = "Bedroom"
$ws.Range("E21").Value = "In Bed"
$ws.Range("F21").Value = 50
$ws.Range("G21").Value = 2
$ws.Range("H21").Value = "Occupied"

$ws.Range("A22").Value = "'2026-01-28"
$ws.Range("A22").Style = "Normal"
$ws.Range("B22").Value = "17:50:28"
$ws.Range("C22").Value = "17:00"
$ws.Range("D22").Value = "Bedroom"
$ws.Range("E22").Value = "In Bed"
$ws.Range("F22").Value = 120
$ws.Range("G22").Value = 72
$ws.Range("H22").Value = "Occupied"

$ws.Range("A23").Value = "'2026-01-28"
$ws.Range("A23").Style = "Normal"
$ws.Range("B23").Value = "17:50:29"
$ws.Range("C23").Value = "17:00"
$ws.Range("D23").Value = "Bedroom"
$ws.Range("E23").Value = "In Bed"
$ws.Range("F23").Value = 76
$ws.Range("G23").Value = 28
$ws.Range("H23").Value = "Occupied"

$ws.Range("A24").Value = "'2026-01-28"
$ws.Range("A24").Style = "Normal"
$ws.Range("B24").Value = "17:50:30"
$ws.Range("C24").Value = "17:00"
$ws.Range("D24").Value = "Bedroom"
$ws.Range("E24").Value = "In Bed"
$ws.Range("F24").Value = 50
$ws.Range("G24").Value = 2
$ws.Range("H24").Value = "Occupied"

$ws.Range("A25").Value = "'2026-01-28"
$ws.Range("A25").Style = "Normal"
$ws.Range("B25").Value = "17:50:33"
$ws.Range("C25").Value = "17:00"
$ws.Range("D25").Value = "Bedroom"
$ws.Range("E25").Value = "In Bed"
$ws.Range("F25").Value = 49
$ws.Range("G25").Value = 1
$ws.Range("H25").Value = "Occupied"

